$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(13, 4).Value = 'id="Par1">Children were relatively spared during COVID-19 pandemic.

 However, the recently reported hyperinflammatory syndrome with overlapping features of Kawasaki disease and toxic shock syndrome—“Paediatric Inflammatory Multisystem Syndrome-temporally associated with SARS-CoV-2” (PIMS-TS) has caused concern.

 We describe cardiac findings and short-term outcomes in children with PIMS-TS at a tertiary children’s hospital.

 Single-center observational study of children with PIMS-TS from 10th April to 9th May 2020. Data on ECG and echocardiogram were retrospectively analyzed along with demographics, clinical features and blood parameters.

 Fifteen children with median age of 8.8 (IQR 6.4–11.2) years were included, all were from African/Afro-Caribbean, South Asian, Mixed or other minority ethnic groups.

 All showed raised inflammatory/cardiac markers (CRP, ferritin, Troponin I, CK and pro-BNP).

 Transient valve regurgitation was present in 10 patients (67%).

 Left Ventricular ejection fraction was reduced in 12 (80%), fractional shortening in 8 (53%) with resolution in all but 2. Fourteen (93%) had coronary artery abnormalities, with normalization in 6. ECG abnormalities were present in 9 (60%) which normalized in 6 by discharge.

 Ten (67%) needed inotropes and/or vasopressors.

 None needed extracorporeal life support.

 Improvement in cardiac biochemical markers was closely followed by improvement in ECG/echocardiogram.

 All patients were discharged alive and twelve (80%) have been reviewed since.

 Our entire cohort with PIMS-TS had cardiac involvement and this degree of involvement is significantly more than other published series and emphasizes the need for specialist cardiac review.

 We believe that our multi-disciplinary team approach was crucial for the good short-term outcomes.

'
$ws.Cells.Item(13, 5).Value = '[Tristan%Ramcharan%NULL%1,    Oscar%Nolan%NULL%1,    Chui Yi%Lai%NULL%1,    Nanda%Prabhu%NULL%1,    Raghu%Krishnamurthy%NULL%1,    Alex G.%Richter%NULL%1,    Deepthi%Jyothish%NULL%1,    Hari Krishnan%Kanthimathinathan%NULL%2,    Steven B.%Welch%NULL%1,    Scott%Hackett%NULL%1,    Eslam%Al-Abadi%NULL%1,    Barnaby R.%Scholefield%NULL%1,    Ashish%Chikermane%a.chikermane@nhs.net%1]'
$ws.Cells.Item(14, 4).Value = 'Background
id="Par1">A recent increase in children admitted with hypotensive shock and fever in the context of the COVID-19 outbreak requires an urgent characterization and assessment of the involvement of SARS-CoV-2 infection.

 This is a case series performed at 4 academic tertiary care centers in Paris of all the children admitted to the pediatric intensive care unit (PICU) with shock, fever and suspected SARS-CoV-2 infection between April 15th and April 27th, 2020.
Results
id="Par2">20 critically ill children admitted for shock had an acute myocarditis (left ventricular ejection fraction, 35% (25–55); troponin, 269 ng/mL (31–4607)), and arterial hypotension with mainly vasoplegic clinical presentation.

 The first symptoms before PICU admission were intense abdominal pain and fever for 6 days (1–10).

 All children had highly elevated C-reactive protein (&gt; 94 mg/L) and procalcitonin (&gt; 1.6 ng/mL) without microbial cause.

 At least one feature of Kawasaki disease was found in all children (fever, n = 20, skin rash, n = 10; conjunctivitis, n = 6; cheilitis, n = 5; adenitis, n = 2), but none had the typical form.

 SARS-CoV-2 PCR and serology were positive for 10 and 15 children, respectively.

 One child had both negative SARS-CoV-2 PCR and serology, but had a typical SARS-CoV-2 chest tomography scan.

 All children but one needed an inotropic/vasoactive drug support (epinephrine, n = 12; milrinone, n = 10; dobutamine, n = 6, norepinephrine, n = 4) and 8 were intubated.

 All children received intravenous immunoglobulin (2 g per kilogram) with adjuvant corticosteroids (n = 2), IL 1 receptor antagonist (n = 1) or a monoclonal antibody against IL-6 receptor (n = 1).

 All children survived and were afebrile with a full left ventricular function recovery at PICU discharge.


Conclusions
id="Par3">Acute myocarditis with intense systemic inflammation and atypical Kawasaki disease is an emerging severe pediatric disease following SARS-CoV-2 infection.

 Early recognition of this disease is needed and referral to an expert center is recommended.

 A delayed and inappropriate host immunological response is suspected.

 While underlying mechanisms remain unclear, further investigations are required to target an optimal treatment.


'
$ws.Cells.Item(14, 5).Value = '[Marion%Grimaud%NULL%0,    Julie%Starck%NULL%0,    Michael%Levy%NULL%0,    Clémence%Marais%NULL%0,    Judith%Chareyre%NULL%0,    Diala%Khraiche%NULL%0,    Marianne%Leruez-Ville%NULL%0,    Pierre%Quartier%NULL%0,    Pierre Louis%Léger%NULL%0,    Guillaume%Geslain%NULL%0,    Nada%Semaan%NULL%0,    Florence%Moulin%NULL%0,    Matthieu%Bendavid%NULL%0,    Sandrine%Jean%NULL%0,    Géraldine%Poncelet%NULL%0,    Sylvain%Renolleau%NULL%0,    Mehdi%Oualha%mehdi.oualha@aphp.fr%0]'
